$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the "Periodo Mora" (column E) and "Valor Mora" (column F) values
# between row 16 and row 18, leaving row 17 untouched.
$ws.Range("E16").Value = "1905"
$ws.Range("F16").Value = 11042
$ws.Range("E18").Value = "1907"
$ws.Range("F18").Value = 33125
